$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = "Rubber Dome"
$ws.Cells.Item(11, 4).Value = 75
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = 2
$ws.Cells.Item(12, 3).Value = "Membrane"
$ws.Cells.Item(12, 4).Value = 86
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = 3
$ws.Cells.Item(13, 3).Value = "Mechanical"
$ws.Cells.Item(13, 4).Value = 83
$ws.Cells.Item(14, 1).Value = 2
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = "Mechanical"
$ws.Cells.Item(14, 4).Value = 74
$ws.Cells.Item(15, 1).Value = 2
$ws.Cells.Item(15, 2).Value = 2
$ws.Cells.Item(15, 3).Value = "Rubber Dome"
$ws.Cells.Item(15, 4).Value = 85
$ws.Cells.Item(16, 1).Value = 2
$ws.Cells.Item(16, 2).Value = 3
$ws.Cells.Item(16, 3).Value = "Membrane"
$ws.Cells.Item(16, 4).Value = 80
$ws.Cells.Item(17, 1).Value = 3
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 3).Value = "Membrane"
$ws.Cells.Item(17, 4).Value = 97
$ws.Cells.Item(18, 1).Value = 3
$ws.Cells.Item(18, 2).Value = 2
$ws.Cells.Item(18, 3).Value = "Mechanical"
$ws.Cells.Item(18, 4).Value = 98
$ws.Cells.Item(19, 1).Value = 3
$ws.Cells.Item(19, 2).Value = 3
$ws.Cells.Item(19, 3).Value = "Rubber Dome"
$ws.Cells.Item(19, 4).Value = 103
$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(20, 3).Value = "Rubber Dome"
$ws.Cells.Item(20, 4).Value = 84
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = 2
$ws.Cells.Item(21, 3).Value = "Membrane"
$ws.Cells.Item(21, 4).Value = 72
$ws.Cells.Item(22, 1).Value = 1
$ws.Cells.Item(22, 2).Value = 3
$ws.Cells.Item(22, 3).Value = "Mechanical"
$ws.Cells.Item(22, 4).Value = 83
$ws.Cells.Item(23, 1).Value = 2
$ws.Cells.Item(23, 2).Value = 1
$ws.Cells.Item(23, 3).Value = "Mechanical"
$ws.Cells.Item(23, 4).Value = 69
$ws.Cells.Item(24, 1).Value = 2
$ws.Cells.Item(24, 2).Value = 2
$ws.Cells.Item(24, 3).Value = "Rubber Dome"
$ws.Cells.Item(24, 4).Value = 76
$ws.Cells.Item(25, 1).Value = 2
$ws.Cells.Item(25, 2).Value = 3
$ws.Cells.Item(25, 3).Value = "Membrane"
$ws.Cells.Item(25, 4).Value = 87
$ws.Cells.Item(26, 1).Value = 3
$ws.Cells.Item(26, 2).Value = 1
$ws.Cells.Item(26, 3).Value = "Membrane"
$ws.Cells.Item(26, 4).Value = 100
$ws.Cells.Item(27, 1).Value = 3
$ws.Cells.Item(27, 2).Value = 2
$ws.Cells.Item(27, 3).Value = "Mechanical"
$ws.Cells.Item(27, 4).Value = 88
$ws.Cells.Item(28, 1).Value = 3
$ws.Cells.Item(28, 2).Value = 3
$ws.Cells.Item(28, 3).Value = "Rubber Dome"
$ws.Cells.Item(28, 4).Value = 96

# Scroll the view and set the selection to match the new data range
$excel.Goto($ws.Range("D28"), $true)
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
